# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
# Both sheets share the same rows 2,5,9,12,13,15,17,21,22 in column F.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1114
    5  = 8815
    9  = 305
    12 = 24
    13 = 3676
    15 = 371
    17 = 3008
    21 = 218
    22 = 2466
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
